$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 304, pushing existing rows 304-379 down to 305-380
$ws.Rows("304:304").Insert()

# Populate the newly inserted row 304 with the new record's data
$ws.Range("A304").Value = 4
$ws.Range("B304").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C304").Value = "Los Lagos"
$ws.Range("D304").Value = 44932
$ws.Range("D304").NumberFormat = $ws.Range("D305").NumberFormat
$ws.Range("E304").Value = 10
$ws.Range("F304").Value = 100112037
$ws.Range("G304").Value = "Cebollín"
$ws.Range("H304").Value = "Sin especificar"
$ws.Range("I304").Value = "Primera"
$ws.Range("J304").Value = 140
$ws.Range("K304").Value = 6000
$ws.Range("L304").Value = 6000
$ws.Range("M304").Value = 6000
$ws.Range("N304").Value = '$/paquete 36 unidades'
$ws.Range("O304").Value = "Región Metropolitana"
$ws.Range("P304").Value = 167
$ws.Range("Q304").Value = 36
$ws.Range("R304").Value = "Hortaliza"
